# "this added last report 01-03-25"
# Roll the daily sale & stock report forward from 19.02.2025 to 01.03.2025:
# update the two header captions, the day's sale figures (top table),
# the stock-movement table, and the cash/bank/extra summary figures at
# the bottom, then leave the view scrolled/selected near the bottom of
# the sheet the way the author left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header captions (shared strings) ---
$ws.Range("A1").Value = "Mangrove Communication  01.03.2025"
$ws.Range("A10").Value = "DAILY STOCK                         (01/03/2025) "

# --- Top "sale" table (rows 3-6) ---
$ws.Range("C3").Value = 29731
$ws.Range("D3").Value = 7250
$ws.Range("E3").Value = $null
$ws.Range("F3").Value = $null

$ws.Range("C4").Value = 29653
$ws.Range("D4").Value = 7247
$ws.Range("E4").Value = 5

$ws.Range("C5").Value = 26110
$ws.Range("D5").Value = 3821

$ws.Range("C6").Value = 48165
$ws.Range("D6").Value = 1977

# Row 7 (C7:G7) are SUM() formulas already on the sheet - they recalc
# automatically from the inputs above.

# --- Stock movement table (rows 13-31) ---
$ws.Range("C13").Value = 53105
$ws.Range("D13").Value = $null

$ws.Range("C14").Value = 383944
$ws.Range("D14").Value = 133659
$ws.Range("E14").Value = $null

$ws.Range("C20").Value = 3470
$ws.Range("D20").Value = 730

$ws.Range("C21").Value = 400
$ws.Range("D21").Value = 152

$ws.Range("C22").Value = 360
$ws.Range("D22").Value = 80

$ws.Range("C24").Value = 40

$ws.Range("C26").Value = 18
$ws.Range("D26").Value = $null

# --- Cash / bank / extra summary (H32:H39) ---
$ws.Range("H33").Value = 225270
$ws.Range("H34").Value = 122674
$ws.Range("H37").Value = 150000

# --- View state: scroll near the bottom rows and select D38 ---
$excel.ActiveWindow.Zoom = 100
$ws.Range("D38").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

$wb.Save()
